$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (style) of an existing header cell onto the
# two new header cells, then set their text.
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Update the existing metric values in row 2
$ws.Range("B2").Value = 0.3628170665402599
$ws.Range("C2").Value = 0.9949857699024791
$ws.Range("D2").Value = 0.4982342734600646
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(n_estimators=100))])"

# Populate the new columns for row 2
$ws.Range("G2").Value = 0.1258822953001072
$ws.Range("H2").Value = 0.988
